# Recreates the "Add files via upload" commit: appends 10 new sticker rows
# (rows 6-16) to the "Стикеры" sheet, introduces the new shared strings,
# a third (blue) font + three new cell styles, four leftover paste-image
# placeholder shapes anchored in column C, and updates the sheet view
# (selection) to match the post-edit state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Стикеры")

# -----------------------------------------------------------------------
# 1) Cell values for the new rows 6-16 (C = file_id, D = key_word)
# -----------------------------------------------------------------------
$ws.Range("C6").Value = 'CAACAgIAAxkBAANEYFnLQxjtNi5MTuMghLPi9mJjD3MAAg0EAALPX4sHmuYS8a7yxGQeBA'
$ws.Range("D6").Value = 'Я не хочу брать Иерусалим'
$ws.Range("C7").Value = 'CAACAgIAAxkBAANBYFnKfcWD9t6m_8'
$ws.Range("D7").Value = 'Смешная шутка'
$ws.Range("C8").Value = 'CAACAgIAAxkBAAOFYFse2EbZxvImJ_jaCeqYhMXZzUUAAkMAA4wPBgUBj55LMpPjEB4E'
$ws.Range("D8").Value = 'Привет'
$ws.Range("C9").Value = 'CAACAgIAAxkBAANKYFsYAAE4EvZcktpJ37Vholo6BLUVAAKKAgACVp29Cj5SbosTxUBnHgQ'
$ws.Range("D9").Value = 'Пока'
$ws.Range("C10").Value = 'CAACAgIAAxkBAAOrYGMFPwhqpYr54oXzp1GzBujaKTYAAlcGAALSWogBjt3QY0E8UxgeBA'
$ws.Range("D10").Value = 'Я люблю аниме'
$ws.Range("C11").Value = 'CAACAgIAAxkBAAOrYGMFPwhqpYr54oXzp1GzBujaKTYAAlcGAALSWogBjt3QY0E8UxgeBA'
$ws.Range("D11").Value = 'я люблю аниме'
$ws.Range("C12").Value = 'CAACAgIAAxkBAAOuYGMFjQcaiyjvWMC1g0usqFnHMEsAAmcBAAIQGm0igOKx4pV8RP0eBA'
$ws.Range("D12").Value = 'Я не люблю аниме'
$ws.Range("C13").Value = 'CAACAgIAAxkBAAOuYGMFjQcaiyjvWMC1g0usqFnHMEsAAmcBAAIQGm0igOKx4pV8RP0eBA'
$ws.Range("D13").Value = 'я не люблю аниме'
$ws.Range("C14").Value = 'CAACAgIAAxkBAAOxYGMF_3KpoGhM_oZyUc8wTbyxl3kAAnwBAAIQGm0iWCyFQx6K9ZYeBA'
$ws.Range("D14").Value = 'Я люблю лоли'
$ws.Range("C15").Value = 'CAACAgIAAxkBAAOxYGMF_3KpoGhM_oZyUc8wTbyxl3kAAnwBAAIQGm0iWCyFQx6K9ZYeBA'
$ws.Range("D15").Value = 'я люблю лоли'

# -----------------------------------------------------------------------
# 2) Styles. Column C in this sheet already carries a "Tahoma 10 black"
#    font (cellXfs index 1 / fontId 1) on C2:C5 -- reuse it verbatim via
#    copy/paste-formats so no redundant font entries are produced, then
#    layer the three brand-new xfs (wrap/valign combos + a new blue font)
#    in the same order the target workbook defines them.
# -----------------------------------------------------------------------
$existingFileIdStyle = $ws.Range("C2")

$existingFileIdStyle.Copy()
$ws.Range("C6").PasteSpecial(-4122)
$existingFileIdStyle.Copy()
$ws.Range("C10").PasteSpecial(-4122)
$existingFileIdStyle.Copy()
$ws.Range("C11").PasteSpecial(-4122)

# New style #2 (cellXfs index 2): same font as above, vertical=top wrap=1 -- first used on C12
$existingFileIdStyle.Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C12").WrapText = $true
$ws.Range("C12").VerticalAlignment = -4160

# New style #3 (cellXfs index 3): same font, horizontal=left vertical=top wrap=1 -- first used on C13
$existingFileIdStyle.Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C13").WrapText = $true
$ws.Range("C13").VerticalAlignment = -4160
$ws.Range("C13").HorizontalAlignment = -4131

# Reuse style #3 / #2 for C14 / C15
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C12").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# New style #4 (cellXfs index 4): new blue font (Tahoma 10, RGB 3A6D99), vertical=top wrap=1 -- C16, left blank
$existingFileIdStyle.Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Font.Color = 10054970
$ws.Range("C16").WrapText = $true
$ws.Range("C16").VerticalAlignment = -4160
$ws.Range("C16").Value = ""
$ws.Application.CutCopyMode = $false

# -----------------------------------------------------------------------
# 3) Row heights for rows 12-15 (Excel auto-fit after wrap-text formatting)
# -----------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 14.25
$ws.Rows.Item(13).RowHeight = 15.75
$ws.Rows.Item(14).RowHeight = 16.5
$ws.Rows.Item(15).RowHeight = 15

# -----------------------------------------------------------------------
# 4) Four leftover "paste image" placeholder autoshapes in column C, one
#    per sticker row (12-15), mirroring the blob: URLs left behind when
#    the Telegram sticker thumbnails were pasted in from the clipboard.
# -----------------------------------------------------------------------
$shp = $ws.Shapes.AddShape(1, $ws.Range("C12").Left, $ws.Range("C12").Top, 24, 24)
$shp.Name = 'AutoShape 1'
$shp.AlternativeText = 'blob:https://web.telegram.org/2d92bd88-1e61-4d99-b3a3-d08b3b858537'
$shp.Fill.Visible = $false
$shp.Line.Visible = $false

$shp = $ws.Shapes.AddShape(1, $ws.Range("C13").Left, $ws.Range("C13").Top, 24, 24)
$shp.Name = 'AutoShape 2'
$shp.AlternativeText = 'blob:https://web.telegram.org/2d92bd88-1e61-4d99-b3a3-d08b3b858537'
$shp.Fill.Visible = $false
$shp.Line.Visible = $false

$shp = $ws.Shapes.AddShape(1, $ws.Range("C14").Left, $ws.Range("C14").Top, 24, 24)
$shp.Name = 'AutoShape 3'
$shp.AlternativeText = 'blob:https://web.telegram.org/2d92bd88-1e61-4d99-b3a3-d08b3b858537'
$shp.Fill.Visible = $false
$shp.Line.Visible = $false

$shp = $ws.Shapes.AddShape(1, $ws.Range("C15").Left, $ws.Range("C15").Top, 24, 24)
$shp.Name = 'AutoShape 4'
$shp.AlternativeText = 'blob:https://web.telegram.org/2d92bd88-1e61-4d99-b3a3-d08b3b858537'
$shp.Fill.Visible = $false
$shp.Line.Visible = $false

# -----------------------------------------------------------------------
# 5) Sheet view: scroll/selection left where the edit ended (E15)
# -----------------------------------------------------------------------
$ws.Activate()
$ws.Range("E15").Select()
